# edit.ps1 - apply the "quan ly de thi" bullet split + drop a stray
# lastRenderedPageBreak marker, per the commit:
#   "chỉnh sửa hàm hiển thị danh sách đề thi, bổ sung thông tin lớp và môn học"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: the single bullet about creating exam papers becomes a short
# "Chức năng quản lý đề thi:" heading bullet followed by two more detailed
# sub-bullets (ilvl 1) - one with the reworded explanation, one describing
# the updated exam-list display.
# ---------------------------------------------------------------------------
$anchor1 = "Chức năng tạo đề thi. Đề thi áp dụng cho lớp. Vì mỗi lớp có thể có nhiều hơn 1 môn học nên sau khi chọn lớp sẽ cần phải chọn môn học để áp dụng cho đề thi."
$find1 = $d.Content.Find
$find1.ClearFormatting()
$found1 = $find1.Execute($anchor1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "edit.ps1: could not find the 'Chức năng tạo đề thi' bullet to split"
}
$targetRange = $find1.Parent
$targetRange.Expand(4) | Out-Null   # wdParagraph -> whole paragraph, incl. the 2nd run's tail sentence

$newBulletsXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Chức năng quản lý đề thi:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Chức năng tạo đề thi</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Vì đề</w:t></w:r><w:r><w:t xml:space="preserve"> thi áp dụng cho lớp. </w:t></w:r><w:r><w:t>M</w:t></w:r><w:r><w:t>ỗi lớp có thể có nhiều hơn 1 môn học nên sau khi chọn lớp sẽ cần phải chọn môn học để áp dụng cho đề thi.</w:t></w:r><w:r><w:t xml:space="preserve"> Hay nói cách khác, danh sách môn thi </w:t></w:r><w:r><w:t>sẽ tương ứng</w:t></w:r><w:r><w:t xml:space="preserve"> lớp học.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Chỉnh sửa lại hàm hiển thị danh sách đề thi với thông tin bổ sung: lớp và môn học.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetRange.InsertXML($newBulletsXml)

# ---------------------------------------------------------------------------
# Edit 2: remove the stray <w:lastRenderedPageBreak/> that sat in front of
# the lone <w:tab/> run at the start of the User-zone paragraph. That marker
# carries no text of its own, so we delete the run's single tab character
# (which removes the whole run, page-break marker included) and then
# re-insert a clean tab-only run at the same (now collapsed) position.
# ---------------------------------------------------------------------------
$anchor2 = "Ngoài những chức năng như version trước"
$find2 = $d.Content.Find
$find2.ClearFormatting()
$found2 = $find2.Execute($anchor2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "edit.ps1: could not find the User-zone paragraph with lastRenderedPageBreak"
}
$paraRange = $find2.Parent
$paraRange.Expand(4) | Out-Null   # wdParagraph -> whole paragraph, incl. the leading tab run
$paraStart = $paraRange.Start

$leadChar = $d.Range($paraStart, $paraStart + 1)
if ($leadChar.Text -ne "`t") {
    throw "edit.ps1: unexpected leading character before the User-zone sentence"
}
$leadChar.Delete()

$insertionPoint = $d.Range($paraStart, $paraStart)
$tabOnlyXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($tabOnlyXml)
